# "updating login page object"
#
# The "suites" sheet's B2/B3 cells hold the most-recent test-report links.
# A new test run appended 5 fresh "/target/capital_bank_ui_smoke*.html"
# report links to the shared-string pool; the cell that used to show the
# placeholder "null" (B2) now shows the newest of those 5 links, and the
# cell that used to show the previous newest report link (B3, the old
# "find_orchards" report) now shows "null".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("suites")

# Append the 5 new report links (in order) so the shared-string table
# grows by 5 entries, the last of which is what B2 ends up pointing to.
$ws.Range("B2").Value = "/target/capital_bank_ui_smokeCnVpPnjAEhHSZUYP.html"
$ws.Range("B2").Value = "/target/capital_bank_ui_smokewQkoayDSNehaUcJx.html"
$ws.Range("B2").Value = "/target/capital_bank_ui_smokeKHaoDVQdroGjXngu.html"
$ws.Range("B2").Value = "/target/capital_bank_ui_smokeQHTkfftiDnLtmkFW.html"
$ws.Range("B2").Value = "/target/capital_bank_ui_smokeLfPBgUEGiAPfgYWC.html"

# B3 reverts to the "null" placeholder that B2 used to hold.
$ws.Range("B3").Value = "null"
